# GoodInfo_v2 - 2021-12-27 未完成
# Append a new data row (row 3) to the sheet with the 2021-12-27 entry.
#
# Values are entered with a leading apostrophe so Excel stores them as
# literal text (matching the existing rows, e.g. "2021-12-24", "-0.7%")
# instead of auto-converting "2021-12-27" into a date serial or "0.25%"
# into a percentage number. The Style is reset to "Normal" afterwards so
# the new cells don't pick up an incidental quote-prefix/number-format
# style and stay on the default style, just like the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "'2021-12-27"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = "'0.25%"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").Value = "'"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = "'"
$ws.Range("D3").Style = "Normal"
